$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new quarter: insert two blank columns before column D.
# The existing quarters in D:K shift right to F:M.
$ws.Columns("D:E").Insert()

# Copy number/date formatting from the (now shifted) column F into the
# two new columns so the new quarters match the existing look.
$ws.Range("F5:F102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Populate the two new quarter columns (D = latest quarter, E = prior one).
$ws.Range("D7").Value = 43434
$ws.Range("E7").Value = 43343
$ws.Range("D8").Value = 8484300
$ws.Range("E8").Value = 8755800
$ws.Range("D9").Value = 8013600
$ws.Range("E9").Value = 8416700
$ws.Range("D10").Value = 470700
$ws.Range("E10").Value = 339100
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = -6400
$ws.Range("E14").Value = -18800
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 8169800
$ws.Range("E17").Value = 8583500
$ws.Range("D18").Value = 314500
$ws.Range("E18").Value = 172300
$ws.Range("D20").Value = 91600
$ws.Range("E20").Value = 43500
$ws.Range("D21").Value = 524700
$ws.Range("E21").Value = 335600
$ws.Range("D22").Value = 38900
$ws.Range("E22").Value = 19000
$ws.Range("D23").Value = 367200
$ws.Range("E23").Value = 196700
$ws.Range("D24").Value = 20100
$ws.Range("E24").Value = 18400
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 347100
$ws.Range("E26").Value = 178300
$ws.Range("D27").Value = 347500
$ws.Range("E27").Value = 178200
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = 21600
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -91600
$ws.Range("E32").Value = -43500
$ws.Range("D33").Value = 347500
$ws.Range("E33").Value = 199800
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 347500
$ws.Range("E35").Value = 199800
$ws.Range("D38").Value = 43434
$ws.Range("E38").Value = 43343
$ws.Range("D41").Value = 266200
$ws.Range("E41").Value = 450600
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 2686100
$ws.Range("E43").Value = 2460400
$ws.Range("D44").Value = 3184400
$ws.Range("E44").Value = 2768600
$ws.Range("D45").Value = 1051000
$ws.Range("E45").Value = 1013500
$ws.Range("D46").Value = 7187700
$ws.Range("E46").Value = 6693200
$ws.Range("D47").Value = 3774500
$ws.Range("E47").Value = 3923900
$ws.Range("D48").Value = 5078300
$ws.Range("E48").Value = 5141700
$ws.Range("D49").Value = 28600
$ws.Range("E49").Value = 167800
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 784600
$ws.Range("E52").Value = 454500
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 16853800
$ws.Range("E54").Value = 16381200
$ws.Range("D57").Value = 2202500
$ws.Range("E57").Value = 1844500
$ws.Range("D58").Value = 2569000
$ws.Range("E58").Value = 2439800
$ws.Range("D59").Value = 1482500
$ws.Range("E59").Value = 1649900
$ws.Range("D60").Value = 6253900
$ws.Range("E60").Value = 5934200
$ws.Range("D61").Value = 1740000
$ws.Range("E61").Value = 1762700
$ws.Range("D62").Value = 567800
$ws.Range("E62").Value = 519300
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 8571000
$ws.Range("E66").Value = 8225600
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 2264000
$ws.Range("E70").Value = 2264000
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 1664000
$ws.Range("E72").Value = 1482000
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 6018700
$ws.Range("E76").Value = 5891500
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43434
$ws.Range("E80").Value = 43343
$ws.Range("D81").Value = 347500
$ws.Range("E81").Value = 199800
$ws.Range("D83").Value = 118600
$ws.Range("E83").Value = 119900
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = -94800
$ws.Range("E89").Value = 1219600
$ws.Range("D91").Value = -104800
$ws.Range("E91").Value = -106300
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -154800
$ws.Range("E94").Value = -7700
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 49100
$ws.Range("E100").Value = -1303000
$ws.Range("D101").Value = -1500
$ws.Range("E101").Value = 7800
$ws.Range("D102").Value = -202000
$ws.Range("E102").Value = -83300

# --- A handful of previously reported figures were restated; correct them
# in their shifted (F:M) locations.
$ws.Range("H8").Value = 8031900
$ws.Range("I8").Value = 8054700
$ws.Range("H9").Value = 7711400
$ws.Range("I9").Value = 8000600
$ws.Range("H10").Value = 320500
$ws.Range("I10").Value = 54100
$ws.Range("I14").Value = 44900
$ws.Range("H17").Value = 7848000
$ws.Range("I17").Value = 8197600
$ws.Range("H18").Value = 183900
$ws.Range("I18").Value = -142900
$ws.Range("H20").Value = 64500
$ws.Range("I20").Value = 46600
$ws.Range("H21").Value = 368600
$ws.Range("I21").Value = 21800
$ws.Range("H23").Value = 207800
$ws.Range("I23").Value = -150200
$ws.Range("H24").Value = 20600
$ws.Range("I24").Value = -43300
$ws.Range("H26").Value = 187200
$ws.Range("I26").Value = -106800
$ws.Range("H27").Value = 187600
$ws.Range("I27").Value = -106900
$ws.Range("H32").Value = -64500
$ws.Range("I32").Value = -46600
$ws.Range("H33").Value = 187600
$ws.Range("I33").Value = -106900
$ws.Range("H35").Value = 187600
$ws.Range("I35").Value = -106900
$ws.Range("I43").Value = 1892200
$ws.Range("I44").Value = 2601600
$ws.Range("I45").Value = 956000
$ws.Range("I46").Value = 5631100
$ws.Range("I49").Value = 171800
$ws.Range("I52").Value = 1518300
$ws.Range("I54").Value = 15818900
$ws.Range("I57").Value = 1991300
$ws.Range("I58").Value = 2141500
$ws.Range("I59").Value = 1349700
$ws.Range("I60").Value = 5482500
$ws.Range("I62").Value = 607300
$ws.Range("I66").Value = 8125800
$ws.Range("I72").Value = 1267800
$ws.Range("I76").Value = 5429100
$ws.Range("H81").Value = 187600
$ws.Range("I81").Value = -106900
$ws.Range("G100").Value = 478500
$ws.Range("H100").Value = 285200
$ws.Range("I100").Value = -773700
$ws.Range("G102").Value = -62100
$ws.Range("H102").Value = 71100
